$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '24.371.86'
Set-TextValue 'E2' '  -1.68%  '

Set-TextValue 'D3' '1.649.27'
Set-TextValue 'E3' '  -3.44%  '

Set-TextValue 'E4' '  -0.12%  '

Set-TextValue 'D5' '310.94'
Set-TextValue 'E5' '  -0.01%  '

Set-TextValue 'D6' '1.000'
Set-TextValue 'E6' '  +0.08%  '

Set-TextValue 'D7' '0.3633'
Set-TextValue 'E7' '  -3.35%  '

Set-TextValue 'D8' '46.76'
Set-TextValue 'E8' '  -5.80%  '

Set-TextValue 'D9' '0.3243'
Set-TextValue 'E9' '  -5.91%  '

Set-TextValue 'D10' '1.120'
Set-TextValue 'E10' '  -7.28%  '

Set-TextValue 'D11' '0.07012'
Set-TextValue 'E11' '  -7.12%  '

Set-TextValue 'E12' '  -0.15%  '

Set-TextValue 'D13' '5.931'
Set-TextValue 'E13' '  -5.97%  '

Set-TextValue 'E14' '  -8.27%  '

Set-TextValue 'D15' '6.580'
Set-TextValue 'E15' '  -6.83%  '

Set-TextValue 'D16' '1.652.98'
Set-TextValue 'E16' '  -3.38%  '

Set-TextValue 'D17' '0.00001039'
Set-TextValue 'E17' '  -8.29%  '

Set-TextValue 'D18' '0.06612'
Set-TextValue 'E18' '  -1.62%  '

Set-TextValue 'E19' '  +0.14%  '

Set-TextValue 'D20' '78.13'
Set-TextValue 'E20' '  -7.93%  '

Set-TextValue 'D21' '5.918'
Set-TextValue 'E21' '  -7.44%  '

Set-TextValue 'E22' '  -9.99%  '

Set-TextValue 'D23' '12.45'
Set-TextValue 'E23' '  -6.16%  '

Set-TextValue 'D24' '24.376.29'
Set-TextValue 'E24' '  -1.75%  '

Set-TextValue 'D25' '2.480'
Set-TextValue 'E25' '  +0.95%  '

Set-TextValue 'D26' '2.320'
Set-TextValue 'E26' '  -17.11%  '

Set-TextValue 'D27' '147.33'
Set-TextValue 'E27' '  -3.20%  '

Set-TextValue 'D28' '18.52'
Set-TextValue 'E28' '  -9.46%  '

Set-TextValue 'D29' '1.835.48'
Set-TextValue 'E29' '  -3.40%  '

Set-TextValue 'B30' 'ImmutableX'
Set-TextValue 'C30' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D30' '1.187'
Set-TextValue 'E30' '  -4.42%  '

Set-TextValue 'B31' 'BitcoinCash'
Set-TextValue 'C31' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D31' '123.65'
Set-TextValue 'E31' '  -6.96%  '

Set-TextValue 'D32' '4.075'
Set-TextValue 'E32' '  -3.23%  '

Set-TextValue 'D33' '5.624'
Set-TextValue 'E33' '  -18.78%  '

Set-TextValue 'D34' '0.08406'
Set-TextValue 'E34' '  -4.73%  '

Set-TextValue 'D35' '1.659'
Set-TextValue 'E35' '  -8.68%  '

Set-TextValue 'D36' '12.21'
Set-TextValue 'E36' '  -11.64%  '

Set-TextValue 'D37' '5.150'
Set-TextValue 'E37' '  -8.04%  '

Set-TextValue 'D38' '1.241'
Set-TextValue 'E38' '  -3.11%  '

Set-TextValue 'D39' '0.06003'
Set-TextValue 'E39' '  -10.08%  '

Set-TextValue 'D40' '0.02211'
Set-TextValue 'E40' '  -8.32%  '

Set-TextValue 'D41' '0.2057'
Set-TextValue 'E41' '  -7.93%  '

Set-TextValue 'D42' '8.136'
Set-TextValue 'E42' '  -13.11%  '

Set-TextValue 'D43' '0.9999'
Set-TextValue 'E43' '  +0.14%  '

Set-TextValue 'D44' '0.5873'
Set-TextValue 'E44' '  -8.95%  '

Set-TextValue 'D45' '3.764'
Set-TextValue 'E45' '  -1.38%  '

Set-TextValue 'D46' '12.58'
Set-TextValue 'E46' '  -10.49%  '

Set-TextValue 'D47' '0.5595'
Set-TextValue 'E47' '  -9.33%  '

Set-TextValue 'D48' '122.08'
Set-TextValue 'E48' '  -6.13%  '

Set-TextValue 'D49' '1.941'
Set-TextValue 'E49' '  -9.13%  '

Set-TextValue 'D50' '0.06887'

Set-TextValue 'B51' 'Tezos'
Set-TextValue 'C51' 'https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz'
Set-TextValue 'D51' '1.173'
Set-TextValue 'E51' '  -4.25%  '
